$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row before current row 3 (THIAGO) for account 004361159 / HFR / 92579.84
$ws.Rows.Item(3).Insert()
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "004361159"
$ws.Cells.Item(3, 2).Value = "HFR"
$ws.Cells.Item(3, 3).Value = 92579.84

# 2) Insert a new row before current row 5 (ZULEIKA, after the shift above) for account 004497875 / HENRIQUE / 20999.16
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "004497875"
$ws.Cells.Item(5, 2).Value = "HENRIQUE"
$ws.Cells.Item(5, 3).Value = 20999.16

# 3) Replace the row that used to hold 005348975 / JULIA / 4060 (now shifted down by 2, to row 7)
#    with 002694089 / VITOR / 4987.29
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "002694089"
$ws.Cells.Item(7, 2).Value = "VITOR"
$ws.Cells.Item(7, 3).Value = 4987.29

# 4) Delete the old row for 004361159 / HFR / 100 further down the sheet
#    (originally row 133, now shifted down by 2 due to the two insertions above -> row 135)
$ws.Rows.Item(135).Delete()
